$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the unused columns:
#   F - "Destination Well Id" (always empty)
#   D - "Destination Plate Type" (redundant with plate barcode)
#   B - "Source Well Id" (always empty)
# Delete from right to left so earlier deletions don't shift later indices.
$ws.Columns("F").Delete()
$ws.Columns("D").Delete()
$ws.Columns("B").Delete()

# Rename headers for the columns that took over the old "Col/Row" data
$ws.Range("B1").Value = "Source Well"
$ws.Range("D1").Value = "Destination Well"

# Update destination plate barcode values (jul17_E -> jul17_G series)
for ($r = 2; $r -le 5; $r++) {
  $ws.Cells.Item($r, 3).Value = "ssdest000000141jul17_G"
}
for ($r = 6; $r -le 7; $r++) {
  $ws.Cells.Item($r, 3).Value = "ssdest000000141jul17_384_G"
}

$ws.Range("E9").Select() | Out-Null
